$wb = $excel.ActiveWorkbook

# --- YDS sheet: append newest game play-by-play yardage logs to the long strings ---
$wsYDS = $wb.Worksheets.Item("YDS")
$wsYDS.Range("B2").Value = "4 9 2 1 2 13 7 4 -4 -4 6 4 1 8 3 4 -2 -3 5 3 -5 4 4 2 4 0 -3 2 4 4 -2 5 11 8 8 10 0 4 13 21 2 3 0 1 5 2 5 5 3 8 3 4 2 10 49 7 2 2 3 9 4 8 4 3 2 8 1 5 1 0 0 3 6 2 7 5 4 8 3 0 3 7 5 9 0 1 3 6 6 12 10 6 1 1 5 2 9 5 6 1 1 8 4 4 11 -1 -2 -2 2 1 6 5 8 5 4 1 0 3 1 0 6 5 1 5 2 3 3 0 7 2 10 1 2 10 9 4 0 1 18 1 7 1 6 25 4 1 1 32 0 0 4 7 2 4 2 0 2 4 -4 6 6 3 2 2 3 9 3 1 20 2 -1 12 3 -2 3 1 2 3 1 4 3 -3 6 1 2 3 11 10 2 20 1 0 2 17 2 1 23 5 12 1 11 -2 4 10 -1 5 17 -2 1 5 -1 3 3 7 -1 4 1 0 6 5 5 -6 5 0 0 7 2 5 3 1 3 8 6 1 2 18 3 1 5 9 11 3 10 6 3 7 1 6 23 4 4 0 6 2 7 2 2 3 1 4 11 3 7 0 6 5 7 10 12 2 0 18 8 1 2 5 4 1 -5 2 0 4 6 3 3 1 4 7 15 4 14 1 5 5 2 6 1 5 7 1 -1 3 2 5 0 1 3 11 8 1 2 17 36 3 0 1 8 3 11 6 1 8 -2 4 4 -1 3 43 3 7 9 7 16 6 0 37 5 6 1 2 0 4 3 10 11 6 2 5 0 7 0 3 0 0 9 2 3 2 4 4 2 3 2 2 3 16 3 7 0 2 9 2 4 5 1 3 7 11 10 5 5 4 2 3 7 2 1 11 2 3 6 2 4 5 2 11 3 12 40 11 7 6 1 1 21 9 11 4 7 5 4 9 1 3 3 5 9 4 2 3 5 1 17 6 4 7 2 3 7 2 5 2 4 1 4 3 1 3 4 -2 4 6 1 36 3 4 10 5 3 3 6 15 3 14 1 12 7 0 0 -2 6 5 4 2 8 4 2 5 2 6 10 -1 9 4 7 -2 3 6 0 1 1 2 -1 1 4 5 4 0 -2 3 0 2 25 8 1 5 1 0 0 2 3 10 5 1 1 11 6 5 1 5 2 5 8 1 2 17 6 1 2 2 2 4 1 11 5 15 3 3 6 2 8 2 10 5 3 9 5 9 3 11 1 1 6 6 14 2 -1 3 2 1 1 1 2 7 4 0 2 6 3 2 2 0 0 4 2 12 0 1 7 1 5 9 -7 8 6 1 0 5 0 7 5 3 6 5 3 3 3 3 8 5 2 7 9 5 7 5 -1 1 -2 1 -3 3 11 3 14 6 3 2 5 5 4 6 10 6 3 3 3 2 5 4 2 3 8 1 3 2 5 8 1 8 5 3 4 6 9 2 9 4 4 4 8 7 1 0 3 2 8 1 5 2 1 5 3 1 23 2 10 5 3 12 3 3 0 11 2 1 8 3 -2 1 5 4 2 3 1 3 4 3 2 0 0 2 18 3 -1 6 2 5 -3 3 9 1 7 1 12 4 -2 20 10 1 4 3 13 3 1 3 -1 3 3 11 7 1 1 3 2 15 2 4 2 3 6 5 1 22 1 3 9 2 9 11 9 9 14 2 4 1 3 8 -1 3 1 6 -1 3 2 44 2 4 2 2 3 1 2 5 10 3 1 6 8 0 2 2 6 4 4 6 13 3 0 5 3 2 6 -1 4 2 14 9 1 -2 4 3 1 9 1 -1 28 3 1 17 7 3 1 2 -2 3 4 8 2 0 0 2 2 1 4 -5 4 2 5 3 -1 3 2 9 9 3 -1 25 24 7 0 4 8 6 5 9 23 4 6 8 0 7 7 5 11 7 1 1 3 10 -1 2 4 7 1 3 6 0 8 9 2 16 5 1 0 -2 7 13 -2 2 2 -1 9 3 0 9 3 8 3 1 -1 2 -1 0 0 44 8 -1 4 -3 5 -2 1 1 9 2 2 6 6 3 -1 3 8 0 0 7 3 -1 4 1 1 -1 3 -3 1 0 -5 3 4 3 1 -1 7 7 4 4 0 2 2 10 3 1 2 19 0 2 1 2 4 0 4 -2 -1 2 0 -4 -2 -1 0 3 2 9 2 2 3 1 6 14 3 -1 30 9 0 -3 0 0 27 4 9 1 5 8 5 10 4 -3 2 -1 5 3 5 1 4 8 6 0 4 6 1 8 13 6 2 -1 11 4 13 3 5 -1 -4 2 2 6 4 0 4 5 0"
$wsYDS.Range("C2").Value = "0 3 6 2 3 2 2 -1 1 3 1 2 6 6 -4 0 21 2 0 0 3 6 4 3 11 1 6 8 4 0 3 -2 1 7 -2 -4 3 5 3 6 0 0 6 6 -5 14 0 0 3 1 13 0 20 12 2 3 -3 0 3 4 2 1 4 4 6 -2 13 2 -3 13 2 2 9 5 3 3 19 5 2 0 0 12 1 7 2 1 7 4 5 -1 0 10 7 1 7 2 3 8 -2 3 5 3 2 5 7 4 0 1 36 3 -2 -4 2 9 6 -2 0 3 1 2 4 13 2 2 0 9 2 -5 5 -1 1 3 7 -1 5 -4 2 1 2 0 7 1 3 5 4 2 5 3 5 3 9 -1 5 1 4 5 38 1 4 1 7 1 -1 3 5 1 -2 0 0 4 7 2 0 0 -3 2 6 0 1 7 -8 -4 5 5 3 -2 1 8 4 3 2 -2 9 6 1 2 2 3 2 3 6 7 4 3 4 3 2 3 3 0 5 7 2 13 3 4 -3 1 1 0 5 1 -1 6 3 0 2 11 -3 0 4 2 -5 2 3 7 0 3 11 0 -1 5 23 1 0 2 2 6 4 2 15 -1 7 5 1 2 1 0 2 0 7 13 1 8 0 -7 3 3 5 2 4 9 7 1 1 1 15 82 24 16 -1 6 2 7 5 3 5 3 2 9 19 7 5 1 1 3 5 1 1 13 -5 3 8 4 11 6 7 1 6 3 16 1 3 4 1 2 1 6 3 2 1 12 4 0 8 24 4 4 1 2 5 5 15 2 15 13 7 2 8 1 3 4 -2 0 5 11 7 2 -3 -2 0 4 7 1 0 3 0 5 6 4 6 5 8 6 10 3 2 3 3 0 2 4 6 4 1 5 3 0 0 3 2 0 1 7 2 4 9 0 5 8 2 6 3 7 2 3 7 3 3 10 5 6 4 1 12 1 2 3 1 0 2 3 4 5 -1 1 6 6 3 3 5 3 3 0 3 5 3 2 6 1 1 0 2 5 -1 -1 16 6 0 7 3 5 9 0 2 -2 0 11 0 -1 1 3 2 -5 3 7 -2 4 3 6 1 2 0 1 17 17 2 2 2 0 -1 7 12 0 4 0 1 7 8 1 17 20 1 1 -3 2 4 1 4 2 1 6 6 5 8 3 0 1 8 5 1 -7 13 5 3 3 1 4 1 2 3 3 2 5 6 16 1 9 2 11 2 -2 -2 16 2 1 3 6 0 0 1 8 12 3 -3 1 1 5 12 2 14 3 6 1 4 5 0 1 -1 5 4 -1 -1 -1 4 13 9 -1 2 -1 0 24 6 2 5 6 3 1 0 4 1 5 0 2 2 1 3 4 4 5 1 2 1 0 1 3 -4 -1 2 1 2 4 2 5 4 -1 6 5 7 1 2 1 5 3 2 2 1 1 0 2 3 4 5 6 -5 0 1 2 3 4 -3 4 1 25 4 4 20 4 3 3 2 4 4 18 9 4 13 4 4 4 5 6 0 2 5 6 0 4 -3 7 5 10 6 2 -6 3 -1 7 2 1 14 1 3 24 6 6 9 5 2 8 -1 5 3 3 -3 1 2 4 2 3 -1 4 10 1 11 15 0 -5 9 6 3 2 -2 0 3 2 1 2 33 2 1 3 4 -2 4 58 2 10 2 3 5 3 3 10 3 2 18 3 1 2 -1 7 2 1 1 14 3 7 2 8 5 3 6 2 3 5 2 -1 6 6 3 8 0 3 12 30 2 3 2 6 0 17 2 2 1 4 2 2 4 3 4 4 2 9 4 -3 1 10 6 2 0 4 4 1 5 3 4 0 0 5 -1 2 0 5 1 3 2 -1 6 5 0 5 21 1 4 15 0 1 1 -1 1 2 -1 0 9 2 6 1 3 2 2 4 3 6 4 9 6 2 19 2"
$wsYDS.Range("B3").Value = "14 -2 5 5 1 4 -2 18 10 12 4 17 8 46 6 5 5 4 38 18 7 29 16 10 25 17 6 8 20 16 4 15 9 8 5 9 3 6 9 17 12 18 8 7 15 3 11 -2 21 7 1 8 7 4 7 7 12 8 10 10 13 14 6 18 52 7 3 2 0 18 3 10 9 16 18 19 8 19 25 10 2 20 15 7 4 12 -2 9 3 19 -5 14 10 4 9 7 4 14 14 6 5 4 9 19 12 6 4 15 9 15 41 3 2 9 9 28 6 6 9 9 2 7 4 15 8 14 10 5 17 3 18 7 22 17 4 6 16 6 14 23 -2 5 7 20 0 8 7 9 7 -1 1 1 47 14 4 5 5 11 5 12 4 8 16 5 12 11 13 3 9 20 4 -3 10 9 10 6 14 6 8 11 6 19 14 40 5 2 11 8 7 -6 21 12 9 20 12 3 5 5 10 -5 5 29 9 1 3 12 3 15 1 8 34 7 5 -1 15 11 3 19 12 15 3 8 19 0 23 9 2 16 44 6 15 9 25 5 19 7 11 6 9 2 3 -2 3 8 24 4 20 16 5 14 7 15 7 6 -3 21 0 18 1 8 7 3 7 18 11 18 10 5 10 6 17 9 2 4 6 -6 2 18 5 6 10 2 9 8 2 9 13 11 8 11 1 22 37 12 20 9 1 10 9 2 9 31 20 5 51 21 8 25 24 9 14 8 12 5 12 17 9 9 19 8 19 19 31 7 11 19 1 9 26 3 44 23 7 6 4 4 41 8 5 5 18 7 8 5 5 9 8 10 17 9 27 18 2 4 9 6 11 4 6 5 13 1 4 9 0 17 11 6 17 7 24 0 6 38 10 7 3 11 6 5 10 5 5 14 12 10 9 7 -1 13 4 4 56 10 16 -3 9 7 13 16 0 8 8 6 7 6 14 3 17 2 3 3 8 5 6 4 10 8 55 9 7 4 16 3 7 4 27 3 8 23 9 12 17 11 11 6 13 4 9 7 10 12 7 -3 11 9 6 18 2 11 6 5 6 19 15 58 9 21 15 3 16 8 4 11 72 20 11 10 12 49 14 -4 11 12 6 32 19 13 8 6 23 7 14 5 9 14 6 31 15 -5 6 7 13 8 11 13 6 4 5 18 7 16 5 14 12 14 -1 15 1 38 2 14 5 9 4 16 4 -1 9 2 19 5 14 7 19 11 20 4 15 -1 14 19 9 7 6 5 9 25 5 7 22 8 8 10 8 6 7 6 46 11 3 12 21 21 8 10 34 26 22 5 17 15 5 7 6 18 19 5 10 9 26 21 -6 17 2 5 2 10 1 20 6 5 22 4 7 8 14 13 9 7 10 17 9 24 11 -3 7 6 3 7 6 12 -2 -6 25 27 9 24 9 1 9 30 6 1 5 13 13 70 19 4 2 15 13 6 8 6 18 19 6 21 7 26 5 8 9 7 40 6 7 33 8 4 3 17 5 7 4 17 8 11 12 4 14 -5 6 4 4 56 29 14 7 6 4 7 30 28 18 15 8 8 9 14 9 4 12 2 8 10 18 26 9 34 13 1 3 5 13 5 8 13 10"
$wsYDS.Range("C3").Value = "29 8 6 7 2 10 6 14 9 12 9 23 9 15 5 37 2 3 2 9 22 9 3 6 9 13 13 3 19 15 20 15 5 4 2 12 4 10 31 6 6 6 8 18 2 17 9 5 12 9 10 5 3 48 5 72 16 5 18 14 5 7 3 12 1 23 1 50 7 15 14 31 16 0 7 9 8 1 4 17 1 9 12 17 8 9 49 4 1 6 3 6 7 8 6 14 64 7 8 1 29 10 2 22 74 5 11 19 8 4 8 7 8 3 4 26 2 5 11 7 6 8 5 -1 8 4 9 11 7 50 24 -4 2 9 5 5 16 5 6 7 9 28 11 -3 3 2 6 2 15 17 5 4 7 11 20 7 6 2 2 10 22 0 4 15 2 16 14 6 20 14 12 5 8 4 14 4 14 22 13 8 7 12 11 13 8 7 1 11 -1 -5 27 4 26 22 3 46 10 4 1 17 8 11 17 7 11 8 4 3 8 18 9 13 4 2 12 20 -2 35 3 9 7 31 10 18 24 14 9 10 27 18 22 9 2 7 5 3 39 19 0 6 13 11 6 6 6 6 14 8 4 -5 11 6 5 14 10 23 23 19 1 2 13 4 20 9 5 16 9 9 8 16 7 7 10 16 11 9 7 19 13 25 8 13 8 14 5 1 2 4 6 24 9 4 5 14 21 11 19 14 4 -1 8 14 25 9 8 12 33 9 3 8 44 23 31 22 18 9 28 6 1 11 11 2 21 15 3 8 18 5 5 18 7 6 6 19 5 11 16 3 8 8 14 8 10 10 8 2 19 6 15 14 13 4 6 7 14 4 31 10 19 7 8 5 6 2 7 32 20 20 19 8 5 13 18 19 5 16 7 17 2 10 6 7 13 8 1 4 8 11 11 16 9 -1 14 7 2 9 6 17 31 13 7 3 10 8 27 7 9 4 4 4 4 11 5 22 9 10 4 5 -5 10 13 -4 6 18 4 27 51 19 7 8 17 20 8 4 5 21 1 54 8 5 11 3 28 17 8 18 1 23 20 9 18 8 7 17 13 7 16 10 7 -7 9 3 6 4 6 9 26 27 84 5 -4 5 10 10 1 7 12 9 28 0 8 25 12 7 8 8 4 20 31 4 9 7 44 6 7 16 6 6 41 9 7 -1 8 2 12 8 50 11 34 7 39 13 6 0 8 7 5 16 7 20 7 6 3 12 49 16 13 -4 9 64 3 2 50 9 16 17 10 1 16 16 9 10 -1 2 5 7 2 39 0 8 4 16 14 33 -2 4 14 6 14 14 -1 23 14 9 7 7 26 1 6 12 7 9 4 17 4 7 14 5 21 4 5 6 28 23 9 5 1 41 7 1 21 10 10 5 17 1 4 20 0 12 2 10 1 3 4 13 25 9 4 3 8 3 2 22 16 11 9 1 6 6 27 3 17 14 14 8 13 5 7 7 -2 3 23 5 8 1 5 5 6 14 19 7 5 8 5 16 6 10 3 6 18 11 13 4 4 5 -1 18 7 11 6 7 6 10 9 17 13 40 24 1 5 12 6 17 21 8 -1 10 3 7 17 1 -5 20 5 8 9 3 3 1 26 9 1 3 28 -6 1 16 33 4 3 6 14 13 7 17 1"

# --- ST sheet: append newest game kickoff/punt yardage logs + season totals bump ---
$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B4").Value = "67 63 56 53 67 69 65 69 71 56 62 45 60 66 58 68 67 65 66 68 71 73 57 65 64 66 56 65 60 63 65 63 65 61 62 47 64 63 59 68 69 65 68 59 62 65 60 67 64 69 25 67 65"
$wsST.Range("B5").Value = "30 16 13 17 25 26 15 21 23 11 16 17 8 32 4 21 25 14 18 21 26 42 17 15 22 26 17 23 13 17 19 21 25 27 20 7 18 14 20 22 30 23 21 23 42 17 22 23 18 25 8 36 24"
$wsST.Range("B6").Value = "29 35 26 22 22 30 14 13 29 28 29 18 28 18 14 31 10 0 26 75 21 21 20 20 17 28 27 33 24 17 32 24 39 28 17 33 16 0 16 24 27 33 33 36 24 18 27 16 21 21 13 27 23 26 17 19 22 20 14 15 18 14 26 32 31 9 27 5"
$wsST.Range("D3").Value = "39 45 40 39 50 51 43 30 36 35 46 34 38 28 41 46 46 50 48 50 48 43 31 31 51 45 27 47 44 46 33 36 42 57 37 38 51 49 49 34 37 39 40 49 58 48 56 38 47 51 47 51 47 49 48 42 40 32 51 44 30 45 31 39 55 45 58 44 60 42 52 59 37 51 44 54 56 44 34 55 43 46 60 57 52 53 58 41 41 48 41 56 45 50 63 39 45 38 39 38 41 41 43 62 53 39 54 49 42 52 42 43 53 47 59 52 55 33 62 45 62 48 49 54 45 58 33 46 54 48 46 49 33 49 42 46 45 38 42 54 63 41 31 34 46 52 48"
$wsST.Range("D4").Value = "0 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 9 0 0 12 5 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 7 0 0 0 5 5 0 0 4 0 -2 0 0 -7 -1 1 0 0 8 0 0 0 0 0 0 0 13 0 17 0 10 14 0 0 0 0 8 0 18 12 0 26 0 0 0 0 10 0 8 0 0 0 0 11 0 3 11 0 0 0 0 0 9 0 0 12 10 0 8 0 13 11 10 1 18 7 16 9 0 0 12 0 28 5 0 11 6 11 3 2 9 10 0 0 0 -1 9 0 0 0 0 10 0 0 0 0 0 12 0"
$wsST.Range("D5").Value = "0 18 16 18 7 0 0 14 11 22 0 0 19 15 18 18 33 14 18 0 6 4 0 42 11 14 0 7 8 15 0 0 0 0 2 0 4 14 0 0 19 0 0 0 0 0 14 0 9 0 8 2 4 0 0 17 9 0 0 0 0 0 1 16 10 0 10 7 0 4 7 0 0 11 0 54 7 0 0 0 0 13 0 0 0 0 0 25 0 0 0 0 9 0 0 0 13 14 2 0 2 0 0 0 0 0 12 0 0 17 7 26 0 4 0 0 0 8 0 14 7 12 0 8 19 0 8 0 2 0 11 12 13 0 0 18 0 34 0 0 4 8 0 0 10 0 0 0 0 17 1 21 2 3 0 0 0"
$wsST.Range("B2").Value = 178
$wsST.Range("D2").Value = 147
$wsST.Range("F2").Value = 84
$wsST.Range("G2").Value = 80
$wsST.Range("J2").Value = 49
$wsST.Range("K2").Value = 44
$wsST.Range("B3").Value = 125

# --- OFF sheet: season totals bump for the new game ---
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("C2").Value = 481
$wsOFF.Range("D2").Value = 17
$wsOFF.Range("F2").Value = 159
$wsOFF.Range("G2").Value = 159
$wsOFF.Range("I2").Value = 17
$wsOFF.Range("J2").Value = 82
$wsOFF.Range("N2").Value = 33
$wsOFF.Range("O2").Value = 37
$wsOFF.Range("C3").Value = 352
$wsOFF.Range("F3").Value = 204
$wsOFF.Range("G3").Value = 60
$wsOFF.Range("H3").Value = 62
$wsOFF.Range("I3").Value = 127
$wsOFF.Range("J3").Value = 107
$wsOFF.Range("L3").Value = 488
$wsOFF.Range("M3").Value = 320
$wsOFF.Range("Q3").Value = 1007

# --- DEF sheet: season totals bump for the new game ---
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("C2").Value = 402
$wsDEF.Range("F2").Value = 121
$wsDEF.Range("J2").Value = 68
$wsDEF.Range("N2").Value = 45
$wsDEF.Range("O2").Value = 50
$wsDEF.Range("P2").Value = 26
$wsDEF.Range("B3").Value = 15
$wsDEF.Range("C3").Value = 376
$wsDEF.Range("D3").Value = 13
$wsDEF.Range("E3").Value = 74
$wsDEF.Range("F3").Value = 235
$wsDEF.Range("G3").Value = 68
$wsDEF.Range("H3").Value = 58
$wsDEF.Range("I3").Value = 144
$wsDEF.Range("J3").Value = 129
$wsDEF.Range("L3").Value = 540
$wsDEF.Range("M3").Value = 309
$wsDEF.Range("Q3").Value = 996

# --- TURNS sheet: season totals bump for the new game ---
$wsTURNS = $wb.Worksheets.Item("TURNS")
$wsTURNS.Range("C3").Value = 19
$wsTURNS.Range("D3").Value = 15
$wsTURNS.Range("E3").Value = 15

# --- PEN sheet: season totals bump for the new game ---
$wsPEN = $wb.Worksheets.Item("PEN")
$wsPEN.Range("B2").Value = 29

